# Agrego caso de prueba para una celda que vale 0, que sí debe ser incluida
# en la lista de dicts resultantes: una nueva fila (Yoda, 0) en la hoja
# "Rebeldes".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rebeldes")
$ws.Activate()

$ws.Range("A4").Value = "Yoda"
$ws.Range("B4").Value = 0

# Deja la selección en la celda siguiente a los datos recién agregados.
[void]$ws.Range("A5").Select()
